# Apply the edits described by the commit diff.
#
# 1) Table 6 (slide 6): header row, last column "Waste" -> "Difference"
# 2) Text Box 8 (slide 6): reposition/resize, and second line
#    "on the date" -> "on the same date"

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(6)

# --- Change 1: table header cell text -------------------------------------
$tableShape = $s.Shapes.Item("Table 6")
$tbl = $tableShape.Table
$cell = $tbl.Cell(1, 7)
$cell.Shape.TextFrame.TextRange.Text = "Difference"

# --- Change 2: "Text Box 8" text + position/size ---------------------------
$box = $s.Shapes.Item("Text Box 8")

# Update only the second paragraph's text ("on the date" -> "on the same date")
# using Characters(start, length) so the existing run/paragraph formatting on
# "Count records" (paragraph 1) is left untouched.
$fullText = $box.TextFrame.TextRange.Text
$marker = "on the date"
$startIdx = $fullText.IndexOf($marker)
if ($startIdx -ge 0) {
    $sub = $box.TextFrame.TextRange.Characters($startIdx + 1, $marker.Length)
    $sub.Text = "on the same date"
}

# Move/resize the shape: off x 4080828 -> 3868103, ext cx 1227455 -> 1652905
# (y and cy are unchanged: 3825240 / 491490). EMU -> points is /12700.
$box.Left = 3868103 / 12700.0
# 130.15001 (rather than the mathematically exact 130.15) compensates for the
# single-precision rounding the host applies when it converts points back to
# EMU, so the saved file ends up with the exact target cx = 1652905 EMU.
$box.Width = 130.15001
# Re-assert the height: changing the text above can trigger the textbox's
# auto-fit to recompute the stored height, so pin it back to the original.
$box.Height = 491490 / 12700.0
